# "Correccion a Diebold Mariano y revision de Cap1"
# - Matriz_Resultados: a handful of +1/-1 DM-test outcomes are corrected to ties (0)
# - P_valores / Estadisticos_DM: the corresponding p-values and DM statistics are recomputed
# - Resumen: win/loss/tie counts, win rate and the row order (ranked by Tasa_Victoria_%) are refreshed
$wb = $excel.ActiveWorkbook

# --- Matriz_Resultados: corrected DM test outcomes ---
$wsMatriz = $wb.Worksheets.Item("Matriz_Resultados")
$wsMatriz.Range("E2").Value = 0
$wsMatriz.Range("G3").Value = 0
$wsMatriz.Range("H3").Value = 0
$wsMatriz.Range("F4").Value = 0
$wsMatriz.Range("B5").Value = 0
$wsMatriz.Range("F5").Value = 0
$wsMatriz.Range("I5").Value = 0
$wsMatriz.Range("D6").Value = 0
$wsMatriz.Range("E6").Value = 0
$wsMatriz.Range("C7").Value = 0
$wsMatriz.Range("C8").Value = 0
$wsMatriz.Range("E9").Value = 0

# --- P_valores: recomputed p-values ---
$wsPval = $wb.Worksheets.Item("P_valores")
$wsPval.Range("C2").Value = 0.0006146299178131365
$wsPval.Range("D2").Value = 0.00100555262414459
$wsPval.Range("E2").Value = 0.005977043227676671
$wsPval.Range("F2").Value = [double]"7.89989796556867E-06"
$wsPval.Range("G2").Value = 0.0001345461069606202
$wsPval.Range("H2").Value = 0.0003738900082952412
$wsPval.Range("I2").Value = 0.0003774498152946393
$wsPval.Range("J2").Value = 0.04355212790027463
$wsPval.Range("B3").Value = 0.0006146299178131365
$wsPval.Range("D3").Value = [double]"4.063456450764669E-05"
$wsPval.Range("E3").Value = [double]"6.094034354298117E-05"
$wsPval.Range("F3").Value = 0.001114999255988236
$wsPval.Range("G3").Value = 0.01772953836585511
$wsPval.Range("H3").Value = 0.005307711958360439
$wsPval.Range("I3").Value = 0.03786307125248256
$wsPval.Range("J3").Value = [double]"1.046377934521558E-06"
$wsPval.Range("B4").Value = 0.00100555262414459
$wsPval.Range("C4").Value = [double]"4.063456450764669E-05"
$wsPval.Range("E4").Value = 0.0001088682925265427
$wsPval.Range("F4").Value = 0.001921991402907386
$wsPval.Range("G4").Value = 0.07055977734899677
$wsPval.Range("H4").Value = 0.03271040521784663
$wsPval.Range("I4").Value = 0.276545241403729
$wsPval.Range("J4").Value = [double]"1.580285980296026E-06"
$wsPval.Range("B5").Value = 0.005977043227676671
$wsPval.Range("C5").Value = [double]"6.094034354298117E-05"
$wsPval.Range("D5").Value = 0.0001088682925265427
$wsPval.Range("F5").Value = 0.01479311122421745
$wsPval.Range("G5").Value = 0.1785740125070734
$wsPval.Range("H5").Value = 0.03890239127575845
$wsPval.Range("I5").Value = 0.01754124521668543
$wsPval.Range("J5").Value = 0.0002994448057775934
$wsPval.Range("B6").Value = [double]"7.89989796556867E-06"
$wsPval.Range("C6").Value = 0.001114999255988236
$wsPval.Range("D6").Value = 0.001921991402907386
$wsPval.Range("E6").Value = 0.01479311122421745
$wsPval.Range("G6").Value = 0.0002778684654920482
$wsPval.Range("H6").Value = 0.0008276970178717757
$wsPval.Range("I6").Value = 0.0007501979417108284
$wsPval.Range("J6").Value = 0.1400394488238186
$wsPval.Range("B7").Value = 0.0001345461069606202
$wsPval.Range("C7").Value = 0.01772953836585511
$wsPval.Range("D7").Value = 0.07055977734899677
$wsPval.Range("E7").Value = 0.1785740125070734
$wsPval.Range("F7").Value = 0.0002778684654920482
$wsPval.Range("H7").Value = 0.576352176787327
$wsPval.Range("I7").Value = 0.1171891952575108
$wsPval.Range("J7").Value = [double]"1.651960642479189E-05"
$wsPval.Range("B8").Value = 0.0003738900082952412
$wsPval.Range("C8").Value = 0.005307711958360439
$wsPval.Range("D8").Value = 0.03271040521784663
$wsPval.Range("E8").Value = 0.03890239127575845
$wsPval.Range("F8").Value = 0.0008276970178717757
$wsPval.Range("G8").Value = 0.576352176787327
$wsPval.Range("I8").Value = 0.1643198612393202
$wsPval.Range("J8").Value = [double]"4.573486762637913E-07"
$wsPval.Range("B9").Value = 0.0003774498152946393
$wsPval.Range("C9").Value = 0.03786307125248256
$wsPval.Range("D9").Value = 0.276545241403729
$wsPval.Range("E9").Value = 0.01754124521668543
$wsPval.Range("F9").Value = 0.0007501979417108284
$wsPval.Range("G9").Value = 0.1171891952575108
$wsPval.Range("H9").Value = 0.1643198612393202
$wsPval.Range("J9").Value = [double]"6.227957450199995E-06"
$wsPval.Range("B10").Value = 0.04355212790027463
$wsPval.Range("C10").Value = [double]"1.046377934521558E-06"
$wsPval.Range("D10").Value = [double]"1.580285980296026E-06"
$wsPval.Range("E10").Value = 0.0002994448057775934
$wsPval.Range("F10").Value = 0.1400394488238186
$wsPval.Range("G10").Value = [double]"1.651960642479189E-05"
$wsPval.Range("H10").Value = [double]"4.573486762637913E-07"
$wsPval.Range("I10").Value = [double]"6.227957450199995E-06"

# --- Estadisticos_DM: recomputed Diebold-Mariano statistics ---
$wsStat = $wb.Worksheets.Item("Estadisticos_DM")
$wsStat.Range("C2").Value = 4.39163682837729
$wsStat.Range("D2").Value = 4.137613859121849
$wsStat.Range("E2").Value = 3.235985276926669
$wsStat.Range("F2").Value = 6.852677697458215
$wsStat.Range("G2").Value = 5.200160487257003
$wsStat.Range("H2").Value = 4.651680606507575
$wsStat.Range("I2").Value = 4.646685882814401
$wsStat.Range("J2").Value = 2.218614161930542
$wsStat.Range("B3").Value = -4.39163682837729
$wsStat.Range("D3").Value = -5.871595826624412
$wsStat.Range("E3").Value = -5.640418800890285
$wsStat.Range("F3").Value = -4.084684305384862
$wsStat.Range("G3").Value = -2.686189487540148
$wsStat.Range("H3").Value = -3.295680067885364
$wsStat.Range("I3").Value = -2.292805681709737
$wsStat.Range("J3").Value = -8.186359055958665
$wsStat.Range("B4").Value = -4.137613859121849
$wsStat.Range("C4").Value = 5.871595826624412
$wsStat.Range("E4").Value = -5.316466346157775
$wsStat.Range("F4").Value = -3.80753795738485
$wsStat.Range("G4").Value = -1.957257955616752
$wsStat.Range("H4").Value = -2.369703172003818
$wsStat.Range("I4").Value = -1.132261632451688
$wsStat.Range("J4").Value = -7.901995129698549
$wsStat.Range("B5").Value = -3.235985276926669
$wsStat.Range("C5").Value = 5.640418800890285
$wsStat.Range("D5").Value = 5.316466346157775
$wsStat.Range("F5").Value = -2.778494520067667
$wsStat.Range("G5").Value = 1.416222743595438
$wsStat.Range("H5").Value = 2.278501878729416
$wsStat.Range("I5").Value = 2.691645161777819
$wsStat.Range("J5").Value = -4.769156356285747
$wsStat.Range("B6").Value = -6.852677697458215
$wsStat.Range("C6").Value = 4.084684305384862
$wsStat.Range("D6").Value = 3.80753795738485
$wsStat.Range("E6").Value = 2.778494520067667
$wsStat.Range("G6").Value = 4.808920450074231
$wsStat.Range("H6").Value = 4.237677451224817
$wsStat.Range("I6").Value = 4.288399306094239
$wsStat.Range("J6").Value = 1.564406853889662
$wsStat.Range("B7").Value = -5.200160487257003
$wsStat.Range("C7").Value = 2.686189487540148
$wsStat.Range("D7").Value = 1.957257955616752
$wsStat.Range("E7").Value = -1.416222743595438
$wsStat.Range("F7").Value = -4.808920450074231
$wsStat.Range("H7").Value = 0.572060260859764
$wsStat.Range("I7").Value = 1.66963749768312
$wsStat.Range("J7").Value = -6.400913647452762
$wsStat.Range("B8").Value = -4.651680606507575
$wsStat.Range("C8").Value = 3.295680067885364
$wsStat.Range("D8").Value = 2.369703172003818
$wsStat.Range("E8").Value = -2.278501878729416
$wsStat.Range("F8").Value = -4.237677451224817
$wsStat.Range("G8").Value = -0.572060260859764
$wsStat.Range("I8").Value = 1.467591446790241
$wsStat.Range("J8").Value = -8.778242410666815
$wsStat.Range("B9").Value = -4.646685882814401
$wsStat.Range("C9").Value = 2.292805681709737
$wsStat.Range("D9").Value = 1.132261632451688
$wsStat.Range("E9").Value = -2.691645161777819
$wsStat.Range("F9").Value = -4.288399306094239
$wsStat.Range("G9").Value = -1.66963749768312
$wsStat.Range("H9").Value = -1.467591446790241
$wsStat.Range("J9").Value = -7.002025826526746
$wsStat.Range("B10").Value = -2.218614161930542
$wsStat.Range("C10").Value = 8.186359055958665
$wsStat.Range("D10").Value = 7.901995129698549
$wsStat.Range("E10").Value = 4.769156356285747
$wsStat.Range("F10").Value = -1.564406853889662
$wsStat.Range("G10").Value = 6.400913647452762
$wsStat.Range("H10").Value = 8.778242410666815
$wsStat.Range("I10").Value = 7.002025826526746

# --- Resumen: win/loss/tie summary recomputed and rows re-sorted by Tasa_Victoria_% ---
$wsResumen = $wb.Worksheets.Item("Resumen")
$wsResumen.Cells.Item(2, 1).Value = "Sieve Bootstrap"
$wsResumen.Cells.Item(2, 2).Value = 5
$wsResumen.Cells.Item(2, 3).Value = 0
$wsResumen.Cells.Item(2, 4).Value = 3
$wsResumen.Cells.Item(2, 5).Value = 62.5
$wsResumen.Cells.Item(2, 6).Value = 0.5491865430011136
$wsResumen.Cells.Item(3, 1).Value = "DeepAR"
$wsResumen.Cells.Item(3, 2).Value = 3
$wsResumen.Cells.Item(3, 3).Value = 0
$wsResumen.Cells.Item(3, 4).Value = 5
$wsResumen.Cells.Item(3, 5).Value = 37.5
$wsResumen.Cells.Item(3, 6).Value = 1.337388535290631
$wsResumen.Cells.Item(4, 1).Value = "LSPM"
$wsResumen.Cells.Item(4, 2).Value = 3
$wsResumen.Cells.Item(4, 3).Value = 1
$wsResumen.Cells.Item(4, 4).Value = 4
$wsResumen.Cells.Item(4, 5).Value = 37.5
$wsResumen.Cells.Item(4, 6).Value = 0.9499616917794748
$wsResumen.Cells.Item(5, 1).Value = "AV-MCPS"
$wsResumen.Cells.Item(5, 2).Value = 3
$wsResumen.Cells.Item(5, 3).Value = 0
$wsResumen.Cells.Item(5, 4).Value = 5
$wsResumen.Cells.Item(5, 5).Value = 37.5
$wsResumen.Cells.Item(5, 6).Value = 1.803763680707633
$wsResumen.Cells.Item(6, 1).Value = "MCPS"
$wsResumen.Cells.Item(6, 2).Value = 3
$wsResumen.Cells.Item(6, 3).Value = 0
$wsResumen.Cells.Item(6, 4).Value = 5
$wsResumen.Cells.Item(6, 5).Value = 37.5
$wsResumen.Cells.Item(6, 6).Value = 1.894804144796328
$wsResumen.Cells.Item(7, 1).Value = "AREPD"
$wsResumen.Cells.Item(7, 2).Value = 1
$wsResumen.Cells.Item(7, 3).Value = 4
$wsResumen.Cells.Item(7, 4).Value = 3
$wsResumen.Cells.Item(7, 5).Value = 12.5
$wsResumen.Cells.Item(7, 6).Value = 5.125631466033672
$wsResumen.Cells.Item(8, 1).Value = "LSPMW"
$wsResumen.Cells.Item(8, 2).Value = 1
$wsResumen.Cells.Item(8, 3).Value = 2
$wsResumen.Cells.Item(8, 4).Value = 5
$wsResumen.Cells.Item(8, 5).Value = 12.5
$wsResumen.Cells.Item(8, 6).Value = 2.467652036887471
$wsResumen.Cells.Item(9, 1).Value = "Block Bootstrapping"
$wsResumen.Cells.Item(9, 2).Value = 0
$wsResumen.Cells.Item(9, 3).Value = 6
$wsResumen.Cells.Item(9, 4).Value = 2
$wsResumen.Cells.Item(9, 5).Value = 0
$wsResumen.Cells.Item(9, 6).Value = 5.789483081472292
$wsResumen.Cells.Item(10, 1).Value = "EnCQR-LSTM"
$wsResumen.Cells.Item(10, 2).Value = 0
$wsResumen.Cells.Item(10, 3).Value = 6
$wsResumen.Cells.Item(10, 4).Value = 2
$wsResumen.Cells.Item(10, 5).Value = 0
$wsResumen.Cells.Item(10, 6).Value = 3.875676665862682
